$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1203, 3, 15, 15, 15, 15),
    @(902, 1, 0, 0, 0, 0),
    @(301, 6, 45, 30, 60, 45),
    @(701, 3, 90, 45, 97, 15),
    @(201, 9, 30, 15, 45, 30),
    @(901, 16, 15, 45, 60, 60),
    @(501, 9, 52, 30, 75, 45),
    @(401, 9, 48, 67, 75, 45),
    @(601, 9, 60, 67, 60, 42),
    @(101, 9, 30, 15, 60, 15),
    @(1001, 18, 30, 75, 60, 72),
    @(801, 3, 67, 65, 52, 45),
    @(1201, 2, 10, 10, 10, 10),
    @(1202, 2, 10, 10, 10, 10),
    @(802, 0, 4, 5, 4, 0),
    @(3, 0, 3, 3, 3, 3),
    @(502, 0, 4, 0, 0, 0),
    @(1101, 0, 15, 30, 30, 0),
    @(1, 0, 2, 2, 2, 2),
    @(2, 0, 2, 2, 2, 2)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowData[$c]
    }
}
